# Weekly update: insert the newest "Achicoria" market reading as row 3,
# pushing all prior history down by one row (row 3->4, 4->5, ... 47->48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; existing rows 3..47 shift down to 4..48.
$ws.Rows(3).Insert()

# Populate the new row 3 with the latest reading.
$ws.Cells.Item(3, 1).Value  = 10                              # Mercado ID
$ws.Cells.Item(3, 2).Value  = "Vega Modelo de Temuco"          # Mercado
$ws.Cells.Item(3, 3).Value  = "La Araucanía"                   # Región
$ws.Cells.Item(3, 4).Value  = 44756                            # Fecha
$ws.Cells.Item(3, 5).Value  = 9                                # Codreg
$ws.Cells.Item(3, 6).Value  = 100112010                        # Categoría ID
$ws.Cells.Item(3, 7).Value  = "Achicoria"                      # Categoría
$ws.Cells.Item(3, 8).Value  = "Sin especificar"                # Variedad
$ws.Cells.Item(3, 9).Value  = "Primera"                        # Calidad
$ws.Cells.Item(3, 10).Value = 550                              # Volumen
$ws.Cells.Item(3, 11).Value = 10000                            # Precio mínimo
$ws.Cells.Item(3, 12).Value = 11000                            # Precio máximo
$ws.Cells.Item(3, 13).Value = 10455                            # Precio promedio ponderado
$ws.Cells.Item(3, 14).Value = "$/caja 18 unidades"             # Unidad de comercialización
$ws.Cells.Item(3, 15).Value = "Región Metropolitana"           # Origen
$ws.Cells.Item(3, 16).Value = 581                              # Precio $/Kg
$ws.Cells.Item(3, 17).Value = 18                               # Kg o Unidades
$ws.Cells.Item(3, 18).Value = "Hortaliza"                      # Clasificación
